# Applies the change: remove the "valueQuantity" slice on Observation.value[x]
# within the Elements sheet (row 50), updates Metadata version/date, and fixes
# up the dependent ranges (AutoFilter, FilterDatabase name, ConditionalFormatting).

$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet (version + date) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.1.15-beta"
$meta.Range("B8").Value = "2023-06-07T11:47:17-05:00"

# --- Update Elements sheet: remove the valueQuantity slice row ---
$ws = $wb.Worksheets.Item("Elements")

# Clear the slicing metadata on the parent Observation.value[x] row (row 49),
# since after removing its only slice there is nothing left to describe.
$ws.Range("AB49").Value = ""
$ws.Range("AC49").Value = ""
$ws.Range("AD49").Value = ""
$ws.Range("AE49").Value = ""

# Delete the "Observation.value[x]:valueQuantity" slice header row. This shifts
# every following row up by one (rows 51-90 become rows 50-89).
$ws.Rows.Item(50).Delete()

# The 7 rows that used to be "Observation.value[x]:valueQuantity.*" (now at
# rows 50-56) no longer belong to a slice, so their ID column should match
# their Path column exactly (drop the ":valueQuantity" qualifier).
for ($r = 50; $r -le 56; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 2).Text
}

# --- Fix up ranges that referenced the old AQ90/AI89 extents ---

# Workbook-level defined name used for the sheet's filter database.
$fdb = $wb.Names.Item("Elements!_FilterDatabase")
$fdb.RefersTo = "=Elements!`$A`$1:`$AQ`$89"

# AutoFilter range + column filters (re-apply on the shrunk range).
$ws.AutoFilterMode = $false
$ws.Range("A1:AQ89").AutoFilter(7, "<> ")
$ws.Range("A1:AQ89").AutoFilter(27, @(""), [Microsoft.Office.Interop.Excel.XlAutoFilterOperator]::xlFilterValues)

# Conditional formatting range (was A2:AI89, now A2:AI88).
$cf = $ws.Cells.Item(2, 1).FormatConditions
$cf.Item(1).ModifyAppliesToRange($ws.Range("A2:AI88"))
$cf.Item(2).ModifyAppliesToRange($ws.Range("A2:AI88"))

Write-Host "Done"
